{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = [\n  \"2025-03-20 Thursday\",\n  \"83-14=\",\n  \"58+17=\",\n  \"63-5=\",\n  \"7+86=\",\n  \"8+13=\",\n  \"84-39=\",\n  \"40-12=\",\n  \"63+8=\",\n  \"36+35=\",\n  \"67-38=\",\n  \"93-88=\",\n  \"18+27=\",\n  \"14+39=\",\n  \"20-16=\",\n  \"76-9=\",\n  \"80-49=\",\n  \"14+29=\",\n  \"61-9=\",\n  \"67+19=\",\n  \"76-47=\",\n  \"48+19=\",\n  \"78+4=\",\n  \"7+84=\",\n  \"52-46=\",\n  \"45+26=\",\n  \"83-78=\",\n  \"75-66=\",\n  \"73-49=\",\n  \"92-44=\",\n  \"54+18=\",\n  \"75-48=\",\n  \"18+44=\",\n  \"76-59=\",\n  \"38+26=\",\n  \"58-39=\",\n  \"86-38=\",\n  \"52-19=\",\n  \"6+39=\",\n  \"62-45=\",\n  \"26+39=\",\n  \"50-24=\",\n  \"3+59=\",\n  \"18+49=\",\n  \"7+89=\",\n  \"38+15=\",\n  \"3+59=\",\n  \"57+37=\",\n  \"59+4=\",\n  \"80-17=\",\n  \"77-69=\",\n  \"45+9=\",\n  \"83-75=\",\n  \"84-28=\",\n  \"9+6=\",\n  \"24+7=\",\n  \"37+39=\",\n  \"37+29=\",\n  \"19+66=\",\n  \"78+4=\",\n  \"27+16=\",\n  \"19+8=\",\n  \"93-57=\",\n  \"19+22=\",\n  \"16+16=\",\n  \"30-28=\",\n  \"72-4=\",\n  \"93-24=\",\n  \"7+69=\",\n  \"50-42=\",\n  \"50-15=\",\n  \"91-89=\",\n  \"19+69=\",\n  \"66+15=\",\n  \"23-17=\",\n  \"28+56=\",\n  \"56-47=\",\n  \"94-28=\",\n  \"86-38=\",\n  \"62-15=\",\n  \"74-29=\",\n  \"73-64=\",\n  \"37+35=\",\n  \"56-28=\",\n  \"45-27=\",\n  \"70-15=\",\n  \"92-13=\",\n  \"31-17=\",\n  \"52-17=\",\n  \"38+24=\",\n  \"9+19=\",\n  \"50-43=\",\n  \"56+5=\",\n  \"91-3=\",\n  \"42+39=\",\n  \"64-35=\",\n  \"91-18=\",\n  \"54+8=\",\n  \"53-6=\",\n  \"72-3=\",\n  \"45-36=\",\n];\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\"Unexpected paragraph count: \" + paragraphs.items.length);\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  paragraphs.items[i].insertText(replacements[i], Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph, outside the table)\n$d.Paragraphs.Item(1).Range.Text = '2025-03-20 Thursday'\n\n# Update every arithmetic-expression cell in the practice table, in row-major order\n$values = @(\n    @('83-14=', '58+17=', '63-5=', '7+86=', '8+13='),\n    @('84-39=', '40-12=', '63+8=', '36+35=', '67-38='),\n    @('93-88=', '18+27=', '14+39=', '20-16=', '76-9='),\n    @('80-49=', '14+29=', '61-9=', '67+19=', '76-47='),\n    @('48+19=', '78+4=', '7+84=', '52-46=', '45+26='),\n    @('83-78=', '75-66=', '73-49=', '92-44=', '54+18='),\n    @('75-48=', '18+44=', '76-59=', '38+26=', '58-39='),\n    @('86-38=', '52-19=', '6+39=', '62-45=', '26+39='),\n    @('50-24=', '3+59=', '18+49=', '7+89=', '38+15='),\n    @('3+59=', '57+37=', '59+4=', '80-17=', '77-69='),\n    @('45+9=', '83-75=', '84-28=', '9+6=', '24+7='),\n    @('37+39=', '37+29=', '19+66=', '78+4=', '27+16='),\n    @('19+8=', '93-57=', '19+22=', '16+16=', '30-28='),\n    @('72-4=', '93-24=', '7+69=', '50-42=', '50-15='),\n    @('91-89=', '19+69=', '66+15=', '23-17=', '28+56='),\n    @('56-47=', '94-28=', '86-38=', '62-15=', '74-29='),\n    @('73-64=', '37+35=', '56-28=', '45-27=', '70-15='),\n    @('92-13=', '31-17=', '52-17=', '38+24=', '9+19='),\n    @('50-43=', '56+5=', '91-3=', '42+39=', '64-35='),\n    @('91-18=', '54+8=', '53-6=', '72-3=', '45-36=')\n)\n\n$t = $d.Tables.Item(1)\nif ($t.Rows.Count -ne $values.Count) {\n    throw \"Unexpected row count: $($t.Rows.Count)\"\n}\n\nfor ($r = 1; $r -le $values.Count; $r++) {\n    if ($t.Columns.Count -ne $values[$r-1].Count) {\n        throw \"Unexpected column count: $($t.Columns.Count)\"\n    }\n    for ($c = 1; $c -le $values[$r-1].Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$r-1][$c-1]\n    }\n}\n"}
